$d = $word.ActiveDocument

$replacements = @(
    @("362×9=", "838×9="),
    @("250×8=", "769×5="),
    @("639×7=", "141×3="),
    @("364×6=", "879×7="),
    @("815×8=", "293×6="),
    @("823×8=", "635×2="),
    @("326×8=", "173×2="),
    @("297×4=", "768×5="),
    @("722×8=", "348×7="),
    @("788×9=", "834×4="),
    @("754×3=", "664×7="),
    @("120×6=", "167×6="),
    @("683×7=", "434×6="),
    @("860×3=", "350×3="),
    @("985×2=", "615×4="),
    @("169×7=", "252×4="),
    @("395×8=", "566×9="),
    @("765×5=", "683×2="),
    @("607×3=", "653×3="),
    @("784×7=", "946×4="),
    @("180×9=", "109×9="),
    @("963×7=", "746×2="),
    @("765×2=", "860×8="),
    @("464×3=", "943×2="),
    @("479×4=", "555×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
